$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.844058765215477
$ws.Range("D2").Value = 8.736837116051161
$ws.Range("E2").Value = 13.02915449505992
$ws.Range("F2").Value = 33.02750896517701
$ws.Range("G2").Value = 35.89756582111963
$ws.Range("H2").Value = 15.79332949549297
$ws.Range("I2").Value = 25.77438479969034
$ws.Range("J2").Value = 9.834544528521278
$ws.Range("L2").Value = 9.993759175289808
$ws.Range("M2").Value = 23.86921922030173
$ws.Range("N2").Value = 17.62447888413172
$ws.Range("O2").Value = 25.06665480334934
$ws.Range("C3").Value = 9.87429085068753
$ws.Range("D3").Value = 8.761697348414412
$ws.Range("E3").Value = 13.08315149208493
$ws.Range("F3").Value = 32.98994453066798
$ws.Range("G3").Value = 35.67223678891654
$ws.Range("H3").Value = 15.80934189931994
$ws.Range("I3").Value = 25.76848845433802
$ws.Range("J3").Value = 9.865625138100018
$ws.Range("L3").Value = 10.02419022330545
$ws.Range("M3").Value = 23.25965987104886
$ws.Range("N3").Value = 17.33153407901387
$ws.Range("O3").Value = 25.04279054237214
$ws.Range("C4").Value = 9.894497025917307
$ws.Range("D4").Value = 8.778022091871206
$ws.Range("E4").Value = 13.11803765197161
$ws.Range("F4").Value = 32.97520943015952
$ws.Range("G4").Value = 35.54435068655263
$ws.Range("H4").Value = 15.82215658180585
$ws.Range("I4").Value = 25.77087485076397
$ws.Range("J4").Value = 9.885639224593291
$ws.Range("L4").Value = 10.0437879153155
$ws.Range("M4").Value = 22.87676099589653
$ws.Range("N4").Value = 17.15104903326928
$ws.Range("O4").Value = 25.03409941806191
$ws.Range("C5").Value = 9.903144164633858
$ws.Range("D5").Value = 8.784941460452746
$ws.Range("E5").Value = 13.1326903725186
$ws.Range("F5").Value = 32.97130335297557
$ws.Range("G5").Value = 35.49491564601338
$ws.Range("H5").Value = 15.82812834910382
$ws.Range("I5").Value = 25.77335842669652
$ws.Range("J5").Value = 9.894029731067775
$ws.Range("L5").Value = 10.05200431851118
$ws.Range("M5").Value = 22.71874793818962
$ws.Range("N5").Value = 17.07743500195658
$ws.Range("O5").Value = 25.03205951695356
$ws.Range("C6").Value = 9.904604947405364
$ws.Range("D6").Value = 8.786106543796866
$ws.Range("E6").Value = 13.13514982397213
$ws.Range("F6").Value = 32.97078158601377
$ws.Range("G6").Value = 35.4868701286831
$ws.Range("H6").Value = 15.82916521875148
$ws.Range("I6").Value = 25.77386209006275
$ws.Range("J6").Value = 9.895437156300369
$ws.Range("L6").Value = 10.05338256702093
$ws.Range("M6").Value = 22.69239629100449
$ws.Range("N6").Value = 17.06521018692935
$ws.Range("O6").Value = 25.03181155484325
$ws.Range("C7").Value = 9.894611972634555
$ws.Range("D7").Value = 8.778114327887746
$ws.Range("E7").Value = 13.11823349587558
$ws.Range("F7").Value = 32.97514825029521
$ws.Range("G7").Value = 35.54367307963571
$ws.Range("H7").Value = 15.82223408446499
$ws.Range("I7").Value = 25.77090222661895
$ws.Range("J7").Value = 9.885751431163907
$ws.Range("L7").Value = 10.04389779176521
$ws.Range("M7").Value = 22.87463772736418
$ws.Range("N7").Value = 17.15005639152542
$ws.Range("O7").Value = 25.03406582381274
$ws.Range("C8").Value = 9.854141476060454
$ws.Range("D8").Value = 8.745188971261955
$ws.Range("E8").Value = 13.04741374440801
$ws.Range("F8").Value = 33.01283087691836
$ws.Range("G8").Value = 35.81773005417257
$ws.Range("H8").Value = 15.798231259424
$ws.Range("I8").Value = 25.77110569027726
$ws.Range("J8").Value = 9.845068385782335
$ws.Range("L8").Value = 10.00406270653408
$ws.Range("M8").Value = 23.66094432118468
$ws.Range("N8").Value = 17.52365272629315
$ws.Range("O8").Value = 25.05719102380449
$ws.Range("C9").Value = 9.787834115759251
$ws.Range("D9").Value = 8.689025767988042
$ws.Range("E9").Value = 12.92223765737862
$ws.Range("F9").Value = 33.1525658296425
$ws.Range("G9").Value = 36.43587832893641
$ws.Range("H9").Value = 15.7748419301489
$ws.Range("I9").Value = 25.81910004050023
$ws.Range("J9").Value = 9.772642608495746
$ws.Range("L9").Value = 9.933160605705568
$ws.Range("M9").Value = 25.12668565930233
$ws.Range("N9").Value = 18.2475359574571
$ws.Range("O9").Value = 25.14967602768956
$ws.Range("C10").Value = 9.747096005706975
$ws.Range("D10").Value = 8.652870328031083
$ws.Range("E10").Value = 12.83856692216571
$ws.Range("F10").Value = 33.29491972461639
$ws.Range("G10").Value = 36.93581680252016
$ws.Range("H10").Value = 15.77209831970058
$ws.Range("I10").Value = 25.88323362928522
$ws.Range("J10").Value = 9.723872584602436
$ws.Range("L10").Value = 9.885425445999802
$ws.Range("M10").Value = 26.14762121130933
$ws.Range("N10").Value = 18.76906311125303
$ws.Range("O10").Value = 25.24606804762422
$ws.Range("C11").Value = 9.730299451296183
$ws.Range("D11").Value = 8.637528207826145
$ws.Range("E11").Value = 12.80229244486571
$ws.Range("F11").Value = 33.36816659101643
$ws.Range("G11").Value = 37.17238225691322
$ws.Range("H11").Value = 15.77398191973428
$ws.Range("I11").Value = 25.9186283028718
$ws.Range("J11").Value = 9.70264128637521
$ws.Range("L11").Value = 9.864646620192966
$ws.Range("M11").Value = 26.59817075928361
$ws.Range("N11").Value = 19.00308649174737
$ws.Range("O11").Value = 25.29600906267951
$ws.Range("C12").Value = 9.724188910703425
$ws.Range("D12").Value = 8.631877249840732
$ws.Range("E12").Value = 12.788812474399
$ws.Range("F12").Value = 33.39710997410911
$ws.Range("G12").Value = 37.26320495757874
$ws.Range("H12").Value = 15.77514471156606
$ws.Range("I12").Value = 25.9329197755264
$ws.Range("J12").Value = 9.69473813457129
$ws.Range("L12").Value = 9.856912179696399
$ws.Range("M12").Value = 26.76665685768641
$ws.Range("N12").Value = 19.09116365522377
$ws.Range("O12").Value = 25.31578717419914
$ws.Range("C13").Value = 9.725493802151004
$ws.Range("D13").Value = 8.633087225793961
$ws.Range("E13").Value = 12.79170423519728
$ws.Range("F13").Value = 33.39082309165776
$ws.Range("G13").Value = 37.24359069720514
$ws.Range("H13").Value = 15.774874307155
$ws.Range("I13").Value = 25.92980245327168
$ws.Range("J13").Value = 9.696434149794495
$ws.Range("L13").Value = 9.858571976898455
$ws.Range("M13").Value = 26.73046693625592
$ws.Range("N13").Value = 19.07222000281633
$ws.Range("O13").Value = 25.31148922751049
$ws.Range("C14").Value = 9.729791721835237
$ws.Range("D14").Value = 8.637060118981152
$ws.Range("E14").Value = 12.80117830511332
$ws.Range("F14").Value = 33.3705236913007
$ws.Range("G14").Value = 37.1798298664904
$ws.Range("H14").Value = 15.77406858003028
$ws.Range("I14").Value = 25.91978630727558
$ws.Range("J14").Value = 9.70198835327238
$ws.Range("L14").Value = 9.864007619783365
$ws.Range("M14").Value = 26.61207555508508
$ws.Range("N14").Value = 19.01034396255031
$ws.Range("O14").Value = 25.29761891146322
$ws.Range("C15").Value = 9.732456884226345
$ws.Range("D15").Value = 8.639514302867509
$ws.Range("E15").Value = 12.80701481537658
$ws.Range("F15").Value = 33.35824637315818
$ws.Range("G15").Value = 37.1409338268214
$ws.Range("H15").Value = 15.7736335600026
$ws.Range("I15").Value = 25.91376660915557
$ws.Range("J15").Value = 9.705408246342119
$ws.Range("L15").Value = 9.867354548541657
$ws.Range("M15").Value = 26.53927669144155
$ws.Range("N15").Value = 18.97237018398888
$ws.Range("O15").Value = 25.28923549022772
$ws.Range("C16").Value = 9.748228655499807
$ws.Range("D16").Value = 8.653895197415414
$ws.Range("E16").Value = 12.84097345947181
$ws.Range("F16").Value = 33.29030250810905
$ws.Range("G16").Value = 36.92053437133039
$ws.Range("H16").Value = 15.77203817549915
$ws.Range("I16").Value = 25.8810451633544
$ws.Range("J16").Value = 9.725279256969012
$ws.Range("L16").Value = 9.886802180168162
$ws.Range("M16").Value = 26.11788594316058
$ws.Range("N16").Value = 18.75369744674846
$ws.Range("O16").Value = 25.24292608499181
$ws.Range("C17").Value = 9.758348896506064
$ws.Range("D17").Value = 8.663000361444499
$ws.Range("E17").Value = 12.86226337965937
$ws.Range("F17").Value = 33.25078580225779
$ws.Range("G17").Value = 36.78761589974913
$ws.Range("H17").Value = 15.77186105892083
$ws.Range("I17").Value = 25.86256064650723
$ws.Range("J17").Value = 9.737713526722246
$ws.Range("L17").Value = 9.898972035763345
$ws.Range("M17").Value = 25.85572421148386
$ws.Range("N17").Value = 18.61866605418301
$ws.Range("O17").Value = 25.21607065628902
$ws.Range("C18").Value = 9.764333114028803
$ws.Range("D18").Value = 8.668341443709657
$ws.Range("E18").Value = 12.87467709143227
$ws.Range("F18").Value = 33.22885686948525
$ws.Range("G18").Value = 36.71202971123709
$ws.Range("H18").Value = 15.77205402046073
$ws.Range("I18").Value = 25.85251472711774
$ws.Range("J18").Value = 9.744955258809652
$ws.Range("L18").Value = 9.906059962052939
$ws.Range("M18").Value = 25.70363641783431
$ws.Range("N18").Value = 18.54070079024731
$ws.Range("O18").Value = 25.20119815315071
$ws.Range("C19").Value = 9.766387305673222
$ws.Range("D19").Value = 8.67016771482327
$ws.Range("E19").Value = 12.87890909012271
$ws.Range("F19").Value = 33.22156992537546
$ws.Range("G19").Value = 36.68658830068124
$ws.Range("H19").Value = 15.77217001265567
$ws.Range("I19").Value = 25.84921414885421
$ws.Range("J19").Value = 9.747422637279209
$ws.Range("L19").Value = 9.908474968792643
$ws.Range("M19").Value = 25.65192315168778
$ws.Range("N19").Value = 18.51425425121173
$ws.Range("O19").Value = 25.1962614386356
$ws.Range("C20").Value = 9.757254673368129
$ws.Range("D20").Value = 8.662020335311887
$ws.Range("E20").Value = 12.85997961661589
$ws.Range("F20").Value = 33.25490972141372
$ws.Range("G20").Value = 36.80167630092758
$ws.Range("H20").Value = 15.77184940410773
$ws.Range("I20").Value = 25.86446775923058
$ws.Range("J20").Value = 9.736380580484266
$ws.Range("L20").Value = 9.897667415058317
$ws.Range("M20").Value = 25.88376723832053
$ws.Range("N20").Value = 18.63307188450876
$ws.Range("O20").Value = 25.21887011750598
$ws.Range("C21").Value = 9.728522531957097
$ws.Range("D21").Value = 8.635888876168446
$ws.Range("E21").Value = 12.79838858829224
$ws.Range("F21").Value = 33.37645349944501
$ws.Range("G21").Value = 37.19852490344328
$ws.Range("H21").Value = 15.77429304978773
$ws.Range("I21").Value = 25.92270423847943
$ws.Range("J21").Value = 9.700353244103718
$ws.Range("L21").Value = 9.862407406116446
$ws.Range("M21").Value = 26.64690869400663
$ws.Range("N21").Value = 19.02853378686873
$ws.Range("O21").Value = 25.30166952311777
$ws.Range("C22").Value = 9.711201491787778
$ws.Range("D22").Value = 8.619735822333171
$ws.Range("E22").Value = 12.75962939941683
$ws.Range("F22").Value = 33.46291324117691
$ws.Range("G22").Value = 37.46508296577113
$ws.Range("H22").Value = 15.7785096546764
$ws.Range("I22").Value = 25.9659391550113
$ws.Range("J22").Value = 9.677603713888812
$ws.Range("L22").Value = 9.840144037925779
$ws.Range("M22").Value = 27.13321399595622
$ws.Range("N22").Value = 19.28379558474168
$ws.Range("O22").Value = 25.36082934911078
$ws.Range("C23").Value = 9.720312606177821
$ws.Range("D23").Value = 8.628272389339507
$ws.Range("E23").Value = 12.78017942135417
$ws.Range("F23").Value = 33.41613055918323
$ws.Range("G23").Value = 37.32218268043164
$ws.Range("H23").Value = 15.77601982399491
$ws.Range("I23").Value = 25.9423927321246
$ws.Range("J23").Value = 9.689672879847414
$ws.Range("L23").Value = 9.85195513283538
$ws.Range("M23").Value = 26.8748433777057
$ws.Range("N23").Value = 19.14787502785337
$ws.Range("O23").Value = 25.32879638091394
$ws.Range("C24").Value = 9.757748854746925
$ws.Range("D24").Value = 8.66246307368783
$ws.Range("E24").Value = 12.86101156428782
$ws.Range("F24").Value = 33.25304283643084
$ws.Range("G24").Value = 36.79531700123185
$ws.Range("H24").Value = 15.77185375490755
$ws.Range("I24").Value = 25.86360374306322
$ws.Range("J24").Value = 9.736982915410975
$ws.Range("L24").Value = 9.89825694959025
$ws.Range("M24").Value = 25.87109323526688
$ws.Range("N24").Value = 18.6265600426357
$ws.Range("O24").Value = 25.21760271404456
$ws.Range("C25").Value = 9.804372352541749
$ws.Range("D25").Value = 8.703321360655993
$ws.Range("E25").Value = 12.95464018715392
$ws.Range("F25").Value = 33.10775429106086
$ws.Range("G25").Value = 36.2603482985898
$ws.Range("H25").Value = 15.77863280978132
$ws.Range("I25").Value = 25.80103177517433
$ws.Range("J25").Value = 9.791452782721722
$ws.Range("L25").Value = 9.951573456086592
$ws.Range("M25").Value = 25.12668565930233
$ws.Range("N25").Value = 18.2475359574571
$ws.Range("O25").Value = 25.14967602768956
